# Weekly update: a new price-report row is inserted as the new row 19
# (dated 2022-05-26 / serial 44707), pushing the previously-existing
# rows 19-80 down one position to 20-81. All cell content other than
# position is preserved by the shift; only the brand-new row needs its
# values written explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 19; Excel shifts rows 19:80 down to 20:81
# and copies formatting (incl. the date style on column D) from the row above.
$ws.Rows("19:19").Insert()

# Populate the newly inserted row with the new weekly price report.
$ws.Range("A19").Value = 6
$ws.Range("B19").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C19").Value = "Metropolitana"
$ws.Range("D19").Value = 44707
$ws.Range("E19").Value = 13
$ws.Range("F19").Value = 100114007
$ws.Range("G19").Value = "Jengibre"
$ws.Range("H19").Value = "Sin especificar"
$ws.Range("I19").Value = "Primera"
$ws.Range("J19").Value = 250
$ws.Range("K19").Value = 12000
$ws.Range("L19").Value = 13000
$ws.Range("M19").Value = 12400
$ws.Range("N19").Value = "`$/caja 13 kilos"
$ws.Range("O19").Value = "Perú"
$ws.Range("P19").Value = 954
$ws.Range("Q19").Value = 13
$ws.Range("R19").Value = "Hortaliza"
